# Auto-update draw results: append the 2025-12-14 Pick 3 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

# Write the new values. Columns A and C look numeric ("2025-12-14",
# "251214") but the sheet stores every column as plain text, so a leading
# apostrophe forces text entry instead of Excel's date/number auto-detect.
$ws.Cells.Item($newRow, 1).Value = "'2025-12-14"
$ws.Cells.Item($newRow, 2).Value = "Pick 3"
$ws.Cells.Item($newRow, 3).Value = "'251214"
$ws.Cells.Item($newRow, 4).Value = "0-8-8"
$ws.Cells.Item($newRow, 5).Value = "2025-12-14T21:38:21.682+04:00"

# The apostrophe entry above tags A/C with a quote-prefix style, which the
# previous row doesn't have. Re-apply the prior row's formatting (style
# only, not its values) so the new row matches the sheet's existing look.
$priorRow = $newRow - 1
$ws.Range($ws.Cells.Item($priorRow, 1), $ws.Cells.Item($priorRow, 5)).Copy() | Out-Null
$ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 5)).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
